# Apply updated activity data (re-sorted innings rows) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of new row number -> original row number, describing where the
# runs/balls/fours/sixes values for that row should come from.
$mapping = @{
    2  = 3
    3  = 5
    4  = 9
    5  = 4
    6  = 2
    7  = 14
    8  = 13
    9  = 10
    10 = 7
    11 = 6
    12 = 12
    13 = 11
    14 = 8
    15 = 16
    16 = 15
}

# Capture original C:F values (runs, balls, fours, sixes) for every row
# before any writes happen, so we don't clobber source data while updating.
$original = @{}
foreach ($row in 2..16) {
    $original[$row] = @(
        $ws.Cells.Item($row, 3).Value2,
        $ws.Cells.Item($row, 4).Value2,
        $ws.Cells.Item($row, 5).Value2,
        $ws.Cells.Item($row, 6).Value2
    )
}

foreach ($row in 2..16) {
    $src = $mapping[$row]
    $vals = $original[$src]
    $ws.Cells.Item($row, 3).Value2 = $vals[0]
    $ws.Cells.Item($row, 4).Value2 = $vals[1]
    $ws.Cells.Item($row, 5).Value2 = $vals[2]
    $ws.Cells.Item($row, 6).Value2 = $vals[3]
}
